$d = $word.ActiveDocument

# Locate the four target paragraphs by their distinctive text content
$paraStats = $null
$paraCorr = $null
$paraCharts = $null
$paraGroups = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Estat*sticas descritivas da base de dados escolhida*") {
        $paraStats = $i
    } elseif ($t -like "Uma matriz de correla*es entre as vari*veis quantitativas.*") {
        $paraCorr = $i
    } elseif ($t -like "Ao menos dois diagramas de barras*") {
        $paraCharts = $i
    } elseif ($t -like "Uma an*lise entre grupos com uso da fun*" ) {
        $paraGroups = $i
    }
}

# --- Paragraph 1: "Estatisticas descritivas..." -> highlight whole paragraph ---
$d.Paragraphs($paraStats).Range.Font.HighlightColorIndex = 7

# --- Paragraph 2: "Uma matriz de correlacoes..." -> highlight whole paragraph ---
$d.Paragraphs($paraCorr).Range.Font.HighlightColorIndex = 7

# --- Paragraph 3: "Ao menos dois diagramas..." ---
# Merge the run ending in "...necessariamente no " with the following
# "ambiente Python (" run (same resulting text, but fused into one run).
$pCharts = $d.Paragraphs($paraCharts)
$fConn = $pCharts.Range.Duplicate
$fConn.Find.ClearFormatting()
$fConn.Find.Execute(" e um diagrama de dispersão, à escolha do estudante. Esses gráficos devem ser feitos necessariamente no ambiente Python (")
$fConn.Text = $fConn.Text + "@"
$tailConn = $d.Range($fConn.End - 1, $fConn.End)
$tailConn.Text = ""

# Merge the " Notebook)" run with the following "." run.
$pCharts = $d.Paragraphs($paraCharts)
$fNb = $pCharts.Range.Duplicate
$fNb.Find.ClearFormatting()
$fNb.Find.Execute(" Notebook).")
$fNb.Text = $fNb.Text + "@"
$tailNb = $d.Range($fNb.End - 1, $fNb.End)
$tailNb.Text = ""

# Highlight the whole paragraph (pPr mark + every run).
$d.Paragraphs($paraCharts).Range.Font.HighlightColorIndex = 7

# --- Paragraph 4: "Uma analise entre grupos..." ---
# Merge ". Exemplos: " + the 3 example questions (with their connecting
# space runs) + " Você deve fazer apenas " into a single run.
$pGroups = $d.Paragraphs($paraGroups)
$fEx = $pGroups.Range.Duplicate
$fEx.Find.ClearFormatting()
$fEx.Find.Execute(". Exemplos: Pessoas com maior escolaridade tendem a ter maior segurança financeira percebida? Existe diferença no conhecimento financeiro médio entre homens e mulheres? Pessoas que têm orçamento familiar (Q3 = 1) se percebem menos estressadas financeiramente? Você deve fazer apenas ")
$fEx.Text = $fEx.Text + "@"
$tailEx = $d.Range($fEx.End - 1, $fEx.End)
$tailEx.Text = ""

# Highlight the whole paragraph (pPr mark + every run, including "uma").
$d.Paragraphs($paraGroups).Range.Font.HighlightColorIndex = 7

"done"
